$wb = $excel.ActiveWorkbook

# Sheet: 展览 (row index matches cell row directly: rows 2-16)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1096
$ws1.Range("F3").Value = 4160
$ws1.Range("F5").Value = 326
$ws1.Range("F8").Value = 35
$ws1.Range("F10").Value = 123
$ws1.Range("F11").Value = 302
$ws1.Range("F12").Value = 231
$ws1.Range("F13").Value = 2897
$ws1.Range("F14").Value = 135
$ws1.Range("F15").Value = 1430
$ws1.Range("G15").Value = 66.90000000000001

# Sheet: 全部类型 (rows shifted by +1 vs 展览 from row 7 onward)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1096
$ws4.Range("F3").Value = 4160
$ws4.Range("F5").Value = 326
$ws4.Range("F9").Value = 35
$ws4.Range("F11").Value = 123
$ws4.Range("F12").Value = 302
$ws4.Range("F13").Value = 231
$ws4.Range("F14").Value = 2897
$ws4.Range("F15").Value = 135
$ws4.Range("F16").Value = 1430
$ws4.Range("G16").Value = 66.90000000000001
